$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename header cells: "<name>_old" -> "<name>_FV2310", "<name>_new" -> "<name>_FV2404"
# ---------------------------------------------------------------------------
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -like "*_old") {
        $cell.Value = $val -replace "_old$", "_FV2310"
    } elseif ($val -like "*_new") {
        $cell.Value = $val -replace "_new$", "_FV2404"
    }
}

# ---------------------------------------------------------------------------
# 2) Turn the range into an Excel Table ("Table1") without disturbing the
#    existing header-row formatting/style index (the engine would otherwise
#    bake the live header format into a brand-new dxf, which the target
#    workbook does not have). We stash the format, strip it, build the
#    table, then paste the original formatting back.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("W1")

$ws.Range("A1").Copy()
$scratch.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$headerRange.ClearFormats()

$dataRange = $ws.Range("A1:U54")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$scratch.Copy()
$headerRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$scratch.Clear()

# ---------------------------------------------------------------------------
# 3) Freeze the header row.
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
